# Reorder the "Comment" column to the end of the header block (columns J:P),
# pushing Buffer, ElectrophoresisCondition, GelType, LadderName,
# SampleLoadingAmount and SampleType one column to the left.
#
# Before:  J=Comment, K=Buffer, L=ElectrophoresisCondition, M=GelType,
#          N=LadderName, O=SampleLoadingAmount, P=SampleType
# After:   J=Buffer, K=ElectrophoresisCondition, L=GelType, M=LadderName,
#          N=SampleLoadingAmount, O=SampleType, P=Comment
#
# Only cells whose content actually changes are touched; columns that were
# (and remain) blank are left alone.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: column headers shift left, Comment moves to the end
$ws.Range("J1").Value = "Buffer"
$ws.Range("K1").Value = "ElectrophoresisCondition"
$ws.Range("L1").Value = "GelType"
$ws.Range("M1").Value = "LadderName"
$ws.Range("N1").Value = "SampleLoadingAmount"
$ws.Range("O1").Value = "SampleType"
$ws.Range("P1").Value = "Comment"

# Row 2: french description for "Comment" moves from J2 to P2
$ws.Range("J2").ClearContents()
$ws.Range("P2").Value = "# Commentaire"

# Row 3: type markers for SampleLoadingAmount/SampleType swap columns (N<->O)
$ws.Range("N3").Value = "#integer"
$ws.Range("O3").Value = "#string"

# Row 4: format hint for "Comment" moves from J4 to P4
$ws.Range("J4").ClearContents()
$ws.Range("P4").Value = "# format: texte libre"
